$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fill in the results of the (re-run) trade scan for this existing trade.
$ws.Range("B2").Value = $false          # Profitable -> FALSE
$ws.Range("E2").Value = 104.06          # SellPrice
$ws.Range("F2").Value = -2.7567517054480914   # Price Change %
$ws.Range("G2").Value = $false          # Holding -> FALSE (position closed)

# Row 3: new trade record added by the latest scan.
$ws.Range("C3").Value = 9724.32         # Principle
